# "GPT was changed to LLM" — the workbook has a single worksheet ("Sheet1")
# with a "Source" column (G) whose data rows (G2:G979) all share the text
# "GPT". Replace that value with "LLM" (this updates the shared-string
# table entry, exactly like Excel's Find & Replace would).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2:G979").Replace("GPT", "LLM")

# The author also scrolled the sheet and left the selection on I976 (near
# the bottom of the data) instead of the previous B965. Reproduce the final
# view/selection state.
$ws.Range("I976").Select()
$excel.ActiveWindow.ScrollRow = 956
$excel.ActiveWindow.ScrollColumn = 1
